$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Assigned to" (column C) for Alex's items: background, walls, ball, paddles
$ws.Range("C3").Value = "Alex"
$ws.Range("C4").Value = "Alex"
$ws.Range("C8").Value = "Alex"
$ws.Range("C9").Value = "Alex"

# New art asset row: logo
$ws.Range("B6").Value = "logo"

# "Assigned to" (column C) for Octavio's items: button, logo, title screen, title
$ws.Range("C5").Value = "Octavio"
$ws.Range("C6").Value = "Octavio"
$ws.Range("C11").Value = "Octavio"
$ws.Range("C12").Value = "Octavio"

# "Assigned to" (column C) for Nick's items: all sound assets + music + technical
$ws.Range("C14").Value = "Nick"
$ws.Range("C15").Value = "Nick"
$ws.Range("C16").Value = "Nick"
$ws.Range("C17").Value = "Nick"
$ws.Range("C18").Value = "Nick"
$ws.Range("C19").Value = "Nick"
$ws.Range("C21").Value = "Nick"
$ws.Range("C22").Value = "Nick"
$ws.Range("C24").Value = "Nick"

# score board assigned to Alex/Nick
$ws.Range("C10").Value = "Alex/Nick"

# "In progress" (column D) marks
$ws.Range("D3").Value = "x"
$ws.Range("D4").Value = "x"
$ws.Range("D6").Value = "x"
$ws.Range("D20").Value = "x"

# New sound asset row: place holder sounds
$ws.Range("B20").Value = "place holder sounds"
$ws.Range("C20").Value = "Nick "

$ws.Range("D5").Select()
